$wb = $excel.ActiveWorkbook

# --- registerValidUserSheet: update test data values ---
$wsValid = $wb.Worksheets.Item("registerValidUserSheet")
$wsValid.Range("A2").Value = "testuserAvengers55"
$wsValid.Range("A3").Value = "Qwerty+123458787"

# --- practiceQuestionSheet: wrap text + taller rows for the updated questions ---
$wsPractice = $wb.Worksheets.Item("practiceQuestionSheet")
$wsPractice.Range("B2").WrapText = $true
$wsPractice.Range("B3").WrapText = $true
$wsPractice.Rows.Item(2).RowHeight = 158.4
$wsPractice.Rows.Item(3).RowHeight = 158.4

# --- move the active tab / view state from practiceQuestionSheet to registerValidUserSheet ---
$wsPractice.Activate()
$wsPractice.Application.ActiveWindow.Zoom = 78
$wsPractice.Range("J4").Select()

$wsValid.Activate()
